$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: Copy cell formatting (number format / font / fill / border) from
# existing, similarly-styled rows down onto the three new diary rows
# (26, 27, 28), which previously held only empty, unstyled placeholder cells.
# ---------------------------------------------------------------------------

# Date column (A) -> same style as the other date cells (m/d/yyyy)
$ws.Range("A24").Copy()
$ws.Range("A26:A28").PasteSpecial(-4122)

# Generic text columns (B,C,D,E,F) -> plain text style used throughout
$ws.Range("B11:F11").Copy()
$ws.Range("B26:F26").PasteSpecial(-4122)
$ws.Range("B27:F27").PasteSpecial(-4122)
$ws.Range("B28:F28").PasteSpecial(-4122)

# Row 26's "Achievements" cell (E26) uses the alternate style seen on E19/E24
$ws.Range("E19").Copy()
$ws.Range("E26").PasteSpecial(-4122)

# Row 27's "Time" cell (B27) uses the alternate style seen on B20
$ws.Range("B20").Copy()
$ws.Range("B27").PasteSpecial(-4122)

# "Your Overall Mood" column (G) styles
$ws.Range("G24").Copy()
$ws.Range("G26").PasteSpecial(-4122)

$ws.Range("G20").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G28").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: Fill in the new diary entries' values.
# ---------------------------------------------------------------------------

# --- Row 26: 2/6/2020 ---
$ws.Range("A26").Value = 43867
$ws.Range("B26").Value = "17:00-19:00 in class"
$ws.Range("C26").Value = "N.A."
$ws.Range("D26").Value = "Follow the lecture professor"
$ws.Range("E26").Value = "Learnt to try to be skeptical, drawing examples and simulating at the same time when reviewing code. Know more about other projects, and others’ dilemmas"
$ws.Range("F26").Value = "Learnt more about what an expert would do in reverse engineering. Also after learning this course for 5 weeks, we are gonna have a midterm next week. Will review the knowledges and practice tools that we learnt before in the next few days and prepare for the midterm."
$ws.Range("G26").Value = "Average"

# --- Row 27: 2/13/2020 ---
$ws.Range("A27").Value = 43874
$ws.Range("B27").Value = "17:00-19:00 in class"
$ws.Range("C27").Value = "N.A."
$ws.Range("D27").Value = "Do well in midterm"
$ws.Range("E27").Value = "finished midterm!!!! big achievement!!!"
$ws.Range("F27").Value = "Midterm is kinda hard to describe, since it has many subjective questions. But it’s also the meaning of this course. There is no certain ways to reverse engineer, what we can do is to be subjective and try our best based on some useful concepts. I hope i did it well. Also about the lecture, it is kinda interesting to think about the big picture of one program. Thinking of it stakeholder, functionality and key developers could let us know more than the program itself, like, the community, and the future."
$ws.Range("G27").Value = "Good"

# --- Row 28: 2/19/2020 ---
$ws.Range("A28").Value = 43880
$ws.Range("B28").Value = "13:00-17:00"
$ws.Range("C28").Value = "Soobin, Marc"
$ws.Range("D28").Value = "Finish homework,3 resubmit homework2"
$ws.Range("E28").Value = "finished homework2, and resubmit our homework3"
$ws.Range("F28").Value = "For our previous homework report, we didn’t realize how the report structure will influence reader’s readability. We change the structure of our homework report and make the content more logical this time."
$ws.Range("G28").Value = "Average"
